# Applies the "Updated symbol list" price/volume/hour refresh to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writes $text into $cell as a literal text value (preserving exact
# formatting such as trailing zeros), without leaving the cell's
# number format changed afterwards.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) '244.50'
Set-TextValue $ws.Cells.Item(2, 7) '15'

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) '24.99'
Set-TextValue $ws.Cells.Item(3, 7) '15'

# Row 4
Set-TextValue $ws.Cells.Item(4, 4) '5.131'
Set-TextValue $ws.Cells.Item(4, 7) '15'

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) '0.05647'
Set-TextValue $ws.Cells.Item(5, 7) '15'

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) '6.510'
Set-TextValue $ws.Cells.Item(6, 7) '15'

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) '2.929'
Set-TextValue $ws.Cells.Item(7, 7) '15'

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) '0.8129'
Set-TextValue $ws.Cells.Item(8, 7) '15'

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) '0.8295'
Set-TextValue $ws.Cells.Item(9, 7) '15'

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) '0.1330'
Set-TextValue $ws.Cells.Item(10, 7) '15'

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) '0.06971'
Set-TextValue $ws.Cells.Item(11, 7) '15'

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) '0.02848'
Set-TextValue $ws.Cells.Item(12, 7) '15'

# Row 13
Set-TextValue $ws.Cells.Item(13, 4) '0.09380'
Set-TextValue $ws.Cells.Item(13, 7) '15'

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) '0.001508'
Set-TextValue $ws.Cells.Item(14, 7) '15'

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) '0.009480'
Set-TextValue $ws.Cells.Item(15, 5) '14OneONEBestin24h'
Set-TextValue $ws.Cells.Item(15, 7) '15'

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) '0.006144'
Set-TextValue $ws.Cells.Item(16, 7) '15'

# Row 17
Set-TextValue $ws.Cells.Item(17, 4) '3.502'
Set-TextValue $ws.Cells.Item(17, 7) '15'

# Row 18
Set-TextValue $ws.Cells.Item(18, 7) '15'

# Row 19
Set-TextValue $ws.Cells.Item(19, 7) '15'

# Row 20
Set-TextValue $ws.Cells.Item(20, 4) '0.03186'
Set-TextValue $ws.Cells.Item(20, 7) '15'

# Row 21
Set-TextValue $ws.Cells.Item(21, 4) '0.1321'
Set-TextValue $ws.Cells.Item(21, 7) '15'

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) '3.754'
Set-TextValue $ws.Cells.Item(22, 7) '15'

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) '0.04674'
Set-TextValue $ws.Cells.Item(23, 7) '15'

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) '0.1361'
Set-TextValue $ws.Cells.Item(24, 7) '15'

# Row 25
Set-TextValue $ws.Cells.Item(25, 7) '15'

# Row 26
Set-TextValue $ws.Cells.Item(26, 4) '0.004245'
Set-TextValue $ws.Cells.Item(26, 7) '15'

# Row 27
Set-TextValue $ws.Cells.Item(27, 4) '0.00009712'
Set-TextValue $ws.Cells.Item(27, 5) '26NitroExNTX'
Set-TextValue $ws.Cells.Item(27, 7) '15'

# Row 28
Set-TextValue $ws.Cells.Item(28, 4) '0.0001968'
Set-TextValue $ws.Cells.Item(28, 7) '15'

# Row 29
Set-TextValue $ws.Cells.Item(29, 7) '15'

# Row 30
Set-TextValue $ws.Cells.Item(30, 7) '15'

# Row 31
Set-TextValue $ws.Cells.Item(31, 7) '15'

# Row 32
Set-TextValue $ws.Cells.Item(32, 7) '15'

# Row 33
Set-TextValue $ws.Cells.Item(33, 7) '15'

# Row 34
Set-TextValue $ws.Cells.Item(34, 7) '15'

# Row 35
Set-TextValue $ws.Cells.Item(35, 7) '15'

# Row 36
Set-TextValue $ws.Cells.Item(36, 7) '15'

# Row 37
Set-TextValue $ws.Cells.Item(37, 7) '15'

# Row 38
Set-TextValue $ws.Cells.Item(38, 7) '15'

# Row 39
Set-TextValue $ws.Cells.Item(39, 7) '15'

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) '0.03611'
Set-TextValue $ws.Cells.Item(40, 7) '15'

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) '0.006278'
Set-TextValue $ws.Cells.Item(41, 7) '15'

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) '0.1046'
Set-TextValue $ws.Cells.Item(42, 7) '15'

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) '0.002720'
Set-TextValue $ws.Cells.Item(43, 7) '15'

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) '0.007403'
Set-TextValue $ws.Cells.Item(44, 7) '15'

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) '0.00005280'
Set-TextValue $ws.Cells.Item(45, 7) '15'

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) '0.00000000751'
Set-TextValue $ws.Cells.Item(46, 7) '15'

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) '0.2002'
Set-TextValue $ws.Cells.Item(47, 5) '46CoinbaseStockTokenCOINWorstin24h'
Set-TextValue $ws.Cells.Item(47, 7) '15'

# Row 48
Set-TextValue $ws.Cells.Item(48, 4) '0.002289'
Set-TextValue $ws.Cells.Item(48, 7) '15'

# Row 49
Set-TextValue $ws.Cells.Item(49, 4) '0.00002102'
Set-TextValue $ws.Cells.Item(49, 7) '15'

# Row 50
Set-TextValue $ws.Cells.Item(50, 4) '0.0002002'
Set-TextValue $ws.Cells.Item(50, 7) '15'

# Row 51
Set-TextValue $ws.Cells.Item(51, 7) '15'

